# Change selection on the "DataFromSeleniumEasyURL" sheet (does not need to stay the active tab).
$wb = $excel.ActiveWorkbook
$seleniumSheet = $wb.Worksheets.Item("DataFromSeleniumEasyURL")
$seleniumSheet.Range("H19").Select() | Out-Null

# Add the new "DDDataFromSeleniumEasyURL" sheet right after "DataFromSeleniumEasyURL"
# (i.e. as the new last sheet), matching sheetId 9 / rId7 in the target workbook.
$headerSheet = $wb.Worksheets.Item("AutoCompleteSampleSheet")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "DDDataFromSeleniumEasyURL"

# Copy the bold/yellow header style used on the other data sheets onto A1 before
# writing values, so the new strings land at the end of the shared string table
# in the same order they appear in the sheet.
$headerSheet.Range("A1:A1").Copy($newSheet.Range("A1"))
$newSheet.Range("A1").Value = "Days"
$newSheet.Range("A2").Value = "Wednesday"
$newSheet.Range("A3").Value = "Sunday"
$newSheet.Range("A4").Value = "Monday"
$newSheet.Range("A5").Value = "Saturday"
$newSheet.Range("A6").Value = "Tuesday"

$newSheet.Columns.Item(1).AutoFit()

# Leave the new sheet's selection on A8, with it as the active tab/sheet.
$newSheet.Range("A8").Select() | Out-Null
